$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "thể loại" / "nhà sản xuất" header row ---
$ws.Range("C1:D1").Clear()

# --- Drop the now-unused category (C) / producer (D) columns for the existing games ---
$ws.Range("C2:D11").Clear()
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(8).AutoFit()

# --- Append the ten new games (11-20) ---
$newGames = @(
    "Far Cry 6",
    "The Sims 4: Deluxe Edition",
    "My Time at Sandrock",
    "The Iron Oath",
    "Jurassic World Evolution 2",
    "Bus Simulator 21",
    "TOGETHER BnB",
    "Alien Shooter 2 - New Era",
    "Nigel's Journey : A Working Day",
    "Lost Wing"
)

$row = 12
$num = 11
foreach ($game in $newGames) {
    $ws.Cells.Item($row, 1).Value = $num
    $ws.Cells.Item($row, 2).Value = $game
    $row++
    $num++
}

# --- New narrow spacer column (E) and image-path column (F) ---
$ws.Columns.Item(5).ColumnWidth = 2.75
$ws.Columns.Item(6).ColumnWidth = 13.75

# --- New image reference for the first game ---
$ws.Range("F2").Value = "dummy/game20.jpg"
$ws.Rows.Item(2).RowHeight = 28.8

$null = $ws.Range("F2").Select()
